$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.405.65'
$ws.Range('E2').Value = '  -1.33%  '
$ws.Range('D3').Value = '2.580.08'
$ws.Range('E3').Value = '  -2.02%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '''582.50'
$ws.Range('E5').Value = '  -2.34%  '
$ws.Range('D6').Value = '''166.25'
$ws.Range('E6').Value = '  -1.01%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  -1.40%  '
$ws.Range('D9').Value = '2.579.67'
$ws.Range('E9').Value = '  -2.02%  '
$ws.Range('E10').Value = '  -3.81%  '
$ws.Range('E11').Value = '  +0.09%  '
$ws.Range('D12').Value = '''0.354'
$ws.Range('E12').Value = '  -2.08%  '
$ws.Range('E13').Value = '  -1.78%  '
$ws.Range('E14').Value = '  -4.15%  '
$ws.Range('E16').Value = '  -3.03%  '
$ws.Range('D17').Value = '66.271.37'
$ws.Range('E17').Value = '  -1.45%  '
$ws.Range('D18').Value = '2.565.11'
$ws.Range('E18').Value = '  -2.49%  '
$ws.Range('E19').Value = '  -6.62%  '
$ws.Range('D20').Value = '''7.68'
$ws.Range('E20').Value = '  -5.25%  '
$ws.Range('D21').Value = '''349.67'
$ws.Range('E21').Value = '  -2.59%  '
$ws.Range('D22').Value = '''4.22'
$ws.Range('E22').Value = '  -3.15%  '
$ws.Range('E23').Value = '  -2.29%  '
$ws.Range('E24').Value = '  -0.13%  '
$ws.Range('E25').Value = '  -4.03%  '
$ws.Range('D26').Value = '''68.80'
$ws.Range('E26').Value = '  -2.46%  '
$ws.Range('D27').Value = '''9.92'
$ws.Range('E27').Value = '  -9.00%  '
$ws.Range('D28').Value = '2.714.17'
$ws.Range('E28').Value = '  -1.99%  '
$ws.Range('E29').Value = '  +0.22%  '
$ws.Range('E30').Value = '  -3.25%  '
$ws.Range('D31').Value = '''529.02'
$ws.Range('E31').Value = '  -4.78%  '
$ws.Range('D32').Value = '''8.07'
$ws.Range('E32').Value = '  +1.67%  '
$ws.Range('E33').Value = '  -3.55%  '
$ws.Range('E34').Value = '  -3.61%  '
$ws.Range('E35').Value = '  -3.96%  '
$ws.Range('E36').Value = '  -0.04%  '
$ws.Range('E37').Value = '  -3.63%  '
$ws.Range('D38').Value = '''156.44'
$ws.Range('E38').Value = '  -0.67%  '
$ws.Range('D39').Value = '''18.71'
$ws.Range('E39').Value = '  -2.54%  '
$ws.Range('D40').Value = '''0.359'
$ws.Range('E40').Value = '  -2.12%  '
$ws.Range('E42').Value = '  -1.69%  '
$ws.Range('D43').Value = '''5.08'
$ws.Range('E43').Value = '  -2.17%  '
$ws.Range('E44').Value = '  -0.01%  '
$ws.Range('E45').Value = '  -2.55%  '
$ws.Range('D46').Value = '0.0₆0285'
$ws.Range('D47').Value = '''148.37'
$ws.Range('E47').Value = '  -2.54%  '
$ws.Range('E48').Value = '  -3.93%  '
$ws.Range('E49').Value = '  -3.36%  '
$ws.Range('E50').Value = '  -2.07%  '
$ws.Range('D51').Value = '''0.0760'
$ws.Range('E51').Value = '  -1.54%  '
